# Horarios Línea 141 - actualización de datos (scrape 09:22:34)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Helper: write a block of rows (array of [A,B,C,D,E] arrays)
# into a worksheet starting at $startRow.
# ---------------------------------------------------------------
function Write-Rows($ws, $startRow, $rows) {
  $r = $startRow
  foreach ($row in $rows) {
    for ($c = 0; $c -lt 5; $c++) {
      $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
  }
}

# =================================================================
# Sheet 1: LP1912
# =================================================================
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 09:22:34"
$ws1.Range("A3").Value = "Total filas: 88"

# Reorder swap in the already-scraped block (rows 38/39 "Linea" column)
$ws1.Range("C38").Value = "15_ABASTO"
$ws1.Range("C39").Value = "11_ETCHEVERRY"

# Reorder swap of rows 49/50 (Hora_Scrap, Linea, Minutos updated)
$ws1.Range("A49").Value = "07:13:03"
$ws1.Range("C49").Value = "23_HERNANDEZ"
$ws1.Range("D49").Value = 99

$ws1.Range("A50").Value = "08:52:40"
$ws1.Range("C50").Value = "215B_EL PATO"
$ws1.Range("D50").Value = 0

# Rows 83-93: two existing rows get new scraped content, the previous
# rows 83-85 shift down to 85-87, and six brand-new rows are appended
# (88-93), all stamped with the new scrape time 09:22:34.
$rows83to93 = @(
  @("09:22:34", "10:21", "23_HERNANDEZ", 59, "LP1912"),
  @("09:22:34", "10:25", "16_SANTA ANA", 63, "LP1912"),
  @("08:38:24", "10:29", "15_ABASTO", 111, "LP1912"),
  @("08:45:31", "10:44", "11X44_ETCHEVERRY", 119, "LP1912"),
  @("08:52:40", "10:46", "15_P INDUSTRIAL", 114, "LP1912"),
  @("09:22:34", "10:53", "27_EL RETIRO", 91, "LP1912"),
  @("09:22:34", "10:57", "10_OLMOS", 95, "LP1912"),
  @("09:22:34", "11:01", "81_EL PELIGRO", 99, "LP1912"),
  @("09:22:34", "11:10", "16_P MOR-SANTA ANA", 108, "LP1912"),
  @("09:22:34", "11:14", "14_ABASTO", 112, "LP1912"),
  @("09:22:34", "11:15", "15X38_ABASTO", 113, "LP1912")
)
Write-Rows $ws1 83 $rows83to93

# =================================================================
# Sheet 2: LP1912-215 (only the timestamp changes)
# =================================================================
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 09:22:34"

# =================================================================
# Sheet 3: 6203-6173
# =================================================================
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 09:22:34"
$ws3.Range("A3").Value = "Total filas: 17"

# Rows 16-22: one brand-new row inserted at the top (row 16, stamped
# with the new scrape time) and the previous rows 16-21 shift down to
# 17-22 unchanged.
$rows16to22 = @(
  @("09:22:34", "09:23", "215A_LA PLATA", 1, "L6173"),
  @("08:28:52", "10:12", "215C_LA PLATA", 104, "L6203"),
  @("08:38:24", "10:13", "215C_LA PLATA", 95, "L6203"),
  @("08:52:40", "10:29", "215B_LP-P MOR-1 Y 57", 97, "L6173"),
  @("08:38:24", "10:30", "215B_LP-P MOR-1 Y 57", 112, "L6173"),
  @("08:52:40", "10:30", "215A_LA PLATA", 98, "L6173"),
  @("08:45:31", "10:31", "215A_LA PLATA", 106, "L6173")
)
Write-Rows $ws3 16 $rows16to22
